# Added Validation for input
# added validation for player ID and submit can now handle default values
#
# This populates the "default_input1" column (column I) for the Tab 1
# rows of Table1 with sample/default values, including replacing the old
# placeholder values ("Quantitative" / "Risk rank 1") that are no longer
# used anywhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5: Recommended position -> example default value (replaces "Quantitative")
$ws.Range("I5").Value = "Defence"

# Row 7: Player synergy -> example default value
$ws.Range("I7").Value = "medium"

# Row 9: Agent -> example default value
$ws.Range("I9").Value = "Jorge"

# Row 10: General style of play -> example default value
$ws.Range("I10").Value = "Team player with a very good physical strength, can support a false nine."

# Row 2: Player ID -> default input / validation placeholder
$ws.Range("I2").Value = "Please enter ID"

# Row 3: Player name -> default input / validation placeholder
$ws.Range("I3").Value = "Please enter Name"

# Row 4: Position -> example default value
$ws.Range("I4").Value = "Defensive Midfield"

# Row 6: Player rating -> example default value, numeric (replaces "Risk rank 1")
$ws.Range("I6").Value = 3

# Row 8: Player age -> example default value, numeric
$ws.Range("I8").Value = 25

# Leave the selection on I4, matching the last cell touched in the edit
$ws.Range("I4").Select()
